$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").EntireRow.Insert()

$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44615
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112031
$ws.Range("G37").Value = "Poroto verde"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 28000
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = 29000
$ws.Range("N37").Value = "$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 1160
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
